$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.975.41'
$ws.Range('D2').Style = $s
$s = $ws.Range('E2').Style
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('E2').Style = $s
$s = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.646.26'
$ws.Range('D3').Style = $s
$s = $ws.Range('E3').Style
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E3').Style = $s
$s = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.46'
$ws.Range('D5').Style = $s
$s = $ws.Range('E5').Style
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('E5').Style = $s
$s = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.51'
$ws.Range('D6').Style = $s
$s = $ws.Range('E6').Style
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('E6').Style = $s
$s = $ws.Range('E7').Style
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E7').Style = $s
$s = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.621'
$ws.Range('D8').Style = $s
$s = $ws.Range('E8').Style
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.88%  '
$ws.Range('E8').Style = $s
$s = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.644.12'
$ws.Range('D9').Style = $s
$s = $ws.Range('E9').Style
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.41%  '
$ws.Range('E9').Style = $s
$s = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.119'
$ws.Range('D10').Style = $s
$s = $ws.Range('E10').Style
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.51%  '
$ws.Range('E10').Style = $s
$s = $ws.Range('E12').Style
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.90%  '
$ws.Range('E12').Style = $s
$s = $ws.Range('E13').Style
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('E13').Style = $s
$s = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.64'
$ws.Range('D14').Style = $s
$s = $ws.Range('E14').Style
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('E14').Style = $s
$s = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.122.63'
$ws.Range('D15').Style = $s
$s = $ws.Range('E15').Style
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.30%  '
$ws.Range('E15').Style = $s
$s = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000185'
$ws.Range('D16').Style = $s
$s = $ws.Range('E16').Style
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.54%  '
$ws.Range('E16').Style = $s
$s = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.902.37'
$ws.Range('D17').Style = $s
$s = $ws.Range('E17').Style
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.63%  '
$ws.Range('E17').Style = $s
$s = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.645.72'
$ws.Range('D18').Style = $s
$s = $ws.Range('E18').Style
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('E18').Style = $s
$s = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.21'
$ws.Range('D19').Style = $s
$s = $ws.Range('E19').Style
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('E19').Style = $s
$s = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.74'
$ws.Range('D20').Style = $s
$s = $ws.Range('E20').Style
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.90%  '
$ws.Range('E20').Style = $s
$s = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.55'
$ws.Range('D21').Style = $s
$s = $ws.Range('E21').Style
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.03%  '
$ws.Range('E21').Style = $s
$s = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '346.15'
$ws.Range('D22').Style = $s
$s = $ws.Range('E22').Style
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('E22').Style = $s
$s = $ws.Range('E23').Style
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('E23').Style = $s
$s = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.01'
$ws.Range('D24').Style = $s
$s = $ws.Range('E24').Style
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.39%  '
$ws.Range('E24').Style = $s
$s = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.86'
$ws.Range('D25').Style = $s
$s = $ws.Range('E25').Style
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +5.62%  '
$ws.Range('E25').Style = $s
$s = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000113'
$ws.Range('D26').Style = $s
$s = $ws.Range('E26').Style
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E26').Style = $s
$s = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '605.63'
$ws.Range('D27').Style = $s
$s = $ws.Range('E27').Style
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.85%  '
$ws.Range('E27').Style = $s
$s = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.29'
$ws.Range('D28').Style = $s
$s = $ws.Range('E28').Style
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.83%  '
$ws.Range('E28').Style = $s
$s = $ws.Range('E29').Style
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.30%  '
$ws.Range('E29').Style = $s
$s = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.18'
$ws.Range('D30').Style = $s
$s = $ws.Range('E30').Style
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.52%  '
$ws.Range('E30').Style = $s
$s = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.162'
$ws.Range('D31').Style = $s
$s = $ws.Range('E31').Style
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E31').Style = $s
$s = $ws.Range('E32').Style
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('E32').Style = $s
$s = $ws.Range('E33').Style
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.77%  '
$ws.Range('E33').Style = $s
$s = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.73'
$ws.Range('D34').Style = $s
$s = $ws.Range('E34').Style
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('E34').Style = $s
$s = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.65'
$ws.Range('D35').Style = $s
$s = $ws.Range('E35').Style
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.94%  '
$ws.Range('E35').Style = $s
$s = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.49'
$ws.Range('D36').Style = $s
$s = $ws.Range('E36').Style
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.14%  '
$ws.Range('E36').Style = $s
$s = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.404'
$ws.Range('D37').Style = $s
$s = $ws.Range('E37').Style
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.24%  '
$ws.Range('E37').Style = $s
$s = $ws.Range('B38').Style
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('B38').Style = $s
$s = $ws.Range('C38').Style
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C38').Style = $s
$s = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.77'
$ws.Range('D38').Style = $s
$s = $ws.Range('E38').Style
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.24%  '
$ws.Range('E38').Style = $s
$s = $ws.Range('B39').Style
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('B39').Style = $s
$s = $ws.Range('C39').Style
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('C39').Style = $s
$s = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = $s
$s = $ws.Range('E39').Style
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E39').Style = $s
$s = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.91'
$ws.Range('D40').Style = $s
$s = $ws.Range('E40').Style
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.66%  '
$ws.Range('E40').Style = $s
$s = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '150.72'
$ws.Range('D41').Style = $s
$s = $ws.Range('E41').Style
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.78%  '
$ws.Range('E41').Style = $s
$s = $ws.Range('B42').Style
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('B42').Style = $s
$s = $ws.Range('C42').Style
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('C42').Style = $s
$s = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.56'
$ws.Range('D42').Style = $s
$s = $ws.Range('E42').Style
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.27%  '
$ws.Range('E42').Style = $s
$s = $ws.Range('B43').Style
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'USDe'
$ws.Range('B43').Style = $s
$s = $ws.Range('C43').Style
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('C43').Style = $s
$s = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').Style = $s
$s = $ws.Range('E43').Style
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('E43').Style = $s
$s = $ws.Range('E44').Style
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('E44').Style = $s
$s = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '163.27'
$ws.Range('D45').Style = $s
$s = $ws.Range('E45').Style
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.89%  '
$ws.Range('E45').Style = $s
$s = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '24.13'
$ws.Range('D46').Style = $s
$s = $ws.Range('E46').Style
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.71%  '
$ws.Range('E46').Style = $s
$s = $ws.Range('E47').Style
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.60%  '
$ws.Range('E47').Style = $s
$s = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0591'
$ws.Range('D48').Style = $s
$s = $ws.Range('E48').Style
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.29%  '
$ws.Range('E48').Style = $s
$s = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.635'
$ws.Range('D49').Style = $s
$s = $ws.Range('E49').Style
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('E49').Style = $s
$s = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.100'
$ws.Range('D50').Style = $s
$s = $ws.Range('E50').Style
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.19%  '
$ws.Range('E50').Style = $s
$s = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0249'
$ws.Range('D51').Style = $s
$s = $ws.Range('E51').Style
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.43%  '
$ws.Range('E51').Style = $s
